# Rewrites the 20 quiz question/answer blocks (paragraphs 1-119, every 6th
# paragraph is a blank separator) to the reordered + updated content. The
# two numbering tables after paragraph 120 are left untouched.
$d = $word.ActiveDocument

# Block 1: Quel réalisateur a dirigé le film "Pulp Fiction" ?
$d.Paragraphs(1).Range.Text = "Quel réalisateur a dirigé le film `"Pulp Fiction`" ?"
$d.Paragraphs(2).Range.Text = "Steven Soderbergh"
$d.Paragraphs(3).Range.Text = "David Fincher"
$d.Paragraphs(4).Range.Text = "Quentin Tarantino"
$d.Paragraphs(5).Range.Text = "les frères Coen"

# Block 2: Dans quel film Leonardo DiCaprio prononce la célèbre réplique « Je suis le roi du monde ! » ?
$d.Paragraphs(7).Range.Text = "Dans quel film Leonardo DiCaprio prononce la célèbre réplique « Je suis le roi du monde ! » ?"
$d.Paragraphs(8).Range.Text = "Inception"
$d.Paragraphs(9).Range.Text = "Le Loup de Wall Street"
$d.Paragraphs(10).Range.Text = "Titanic"
$d.Paragraphs(11).Range.Text = "Aviator"

# Block 3: Quel réalisateur a dirigé le film "Le seigneur des anneaux" ?
$d.Paragraphs(13).Range.Text = "Quel réalisateur a dirigé le film `"Le seigneur des anneaux`" ?"
$d.Paragraphs(14).Range.Text = "Guillermo del Toro"
$d.Paragraphs(15).Range.Text = "Peter Jackson"
$d.Paragraphs(16).Range.Text = "George Lucas"
$d.Paragraphs(17).Range.Text = "Ron Howard"

# Block 4: Quel film d'animation a remporté l'Oscar du meilleur film d'animation en 2021 ?
$d.Paragraphs(19).Range.Text = "Quel film d'animation a remporté l'Oscar du meilleur film d'animation en 2021 ?"
$d.Paragraphs(20).Range.Text = "`"Soul`""
$d.Paragraphs(21).Range.Text = "`"En avant`""
$d.Paragraphs(22).Range.Text = "`"Les Mitchell contre les machines`""
$d.Paragraphs(23).Range.Text = "`"Soul`""

# Block 5: Quel réalisateur est connu pour ses films d'animation du Studio Ghibli, comme « Mon voisin Totoro » et « Le Château ambulant » ?
$d.Paragraphs(25).Range.Text = "Quel réalisateur est connu pour ses films d'animation du Studio Ghibli, comme « Mon voisin Totoro » et « Le Château ambulant » ?"
$d.Paragraphs(26).Range.Text = "Mamoru Hosoda"
$d.Paragraphs(27).Range.Text = "Hayao Miyazaki"
$d.Paragraphs(28).Range.Text = "Makoto Shinkai"
$d.Paragraphs(29).Range.Text = "Satoshi Kon"

# Block 6: Quel film a gagné l'oscar du meilleur film en 1994 ?
$d.Paragraphs(31).Range.Text = "Quel film a gagné l'oscar du meilleur film en 1994 ?"
$d.Paragraphs(32).Range.Text = "Pulp Fiction"
$d.Paragraphs(33).Range.Text = "Forrest Gump"
$d.Paragraphs(34).Range.Text = "Schindler's List"
$d.Paragraphs(35).Range.Text = "Les vestiges du jour"

# Block 7: Quelle actrice a joué le rôle de la princesse Leia dans la saga "Star Wars" ?
$d.Paragraphs(37).Range.Text = "Quelle actrice a joué le rôle de la princesse Leia dans la saga `"Star Wars`" ?"
$d.Paragraphs(38).Range.Text = "Sigourney Weaver"
$d.Paragraphs(39).Range.Text = "Meryl Streep"
$d.Paragraphs(40).Range.Text = "Carrie Fisher"
$d.Paragraphs(41).Range.Text = "Natalie Portman"

# Block 8: Quelle actrice a joué le rôle de Katniss Everdeen dans la saga "Hunger Games" ?
$d.Paragraphs(43).Range.Text = "Quelle actrice a joué le rôle de Katniss Everdeen dans la saga `"Hunger Games`" ?"
$d.Paragraphs(44).Range.Text = "Emma Stone"
$d.Paragraphs(45).Range.Text = "Jennifer Aniston"
$d.Paragraphs(46).Range.Text = "Jennifer Lawrence"
$d.Paragraphs(47).Range.Text = "Scarlett Johansson"

# Block 9: Quel est le nom du personnage principal de la saga "Harry Potter" ?
$d.Paragraphs(49).Range.Text = "Quel est le nom du personnage principal de la saga `"Harry Potter`" ?"
$d.Paragraphs(50).Range.Text = "Ron Weasley"
$d.Paragraphs(51).Range.Text = "Hermione Granger"
$d.Paragraphs(52).Range.Text = "Harry Potter"
$d.Paragraphs(53).Range.Text = "Drago Malefoy"

# Block 10: Quel est le titre du premier film d'animation des studios Pixar ?
$d.Paragraphs(55).Range.Text = "Quel est le titre du premier film d'animation des studios Pixar ?"
$d.Paragraphs(56).Range.Text = "Toy Story"
$d.Paragraphs(57).Range.Text = "1001 pattes"
$d.Paragraphs(58).Range.Text = "Monstres et Cie"
$d.Paragraphs(59).Range.Text = "Le monde de Nemo"

# Block 11: Qui a réalisé le film "Psychose" ?
$d.Paragraphs(61).Range.Text = "Qui a réalisé le film `"Psychose`" ?"
$d.Paragraphs(62).Range.Text = "Orson Welles"
$d.Paragraphs(63).Range.Text = "Fritz Lang"
$d.Paragraphs(64).Range.Text = "Alfred Hitchcock"
$d.Paragraphs(65).Range.Text = "Billy Wilder"

# Block 12: Quel acteur a joué le rôle de Joker dans le film "The Dark Knight" ?
$d.Paragraphs(67).Range.Text = "Quel acteur a joué le rôle de Joker dans le film `"The Dark Knight`" ?"
$d.Paragraphs(68).Range.Text = "Jack Nicholson"
$d.Paragraphs(69).Range.Text = "Heath Ledger"
$d.Paragraphs(70).Range.Text = "Joaquin Phoenix"
$d.Paragraphs(71).Range.Text = "Jared Leto"

# Block 13: Quel acteur a joué le rôle de Travis Bickle dans "Taxi Driver" ?
$d.Paragraphs(73).Range.Text = "Quel acteur a joué le rôle de Travis Bickle dans `"Taxi Driver`" ?"
$d.Paragraphs(74).Range.Text = "Harvey Keitel"
$d.Paragraphs(75).Range.Text = "Robert De Niro"
$d.Paragraphs(76).Range.Text = "Al Pacino"
$d.Paragraphs(77).Range.Text = "Joe Pesci"

# Block 14: Quel film a remporté la palme d'or au festival de cannes 2022 ?
$d.Paragraphs(79).Range.Text = "Quel film a remporté la palme d'or au festival de cannes 2022 ?"
$d.Paragraphs(80).Range.Text = "`"As Bestas`""
$d.Paragraphs(81).Range.Text = "`"Armageddon Time`""
$d.Paragraphs(82).Range.Text = "`"Nostalgia`""
$d.Paragraphs(83).Range.Text = "`"Triangle of Sadness`""

# Block 15: Quel réalisateur est connu pour ses films de science-fiction tels que "Inception" et "Interstellar" ?
$d.Paragraphs(85).Range.Text = "Quel réalisateur est connu pour ses films de science-fiction tels que `"Inception`" et `"Interstellar`" ?"
$d.Paragraphs(86).Range.Text = "James Cameron"
$d.Paragraphs(87).Range.Text = "Ridley Scott"
$d.Paragraphs(88).Range.Text = "Christopher Nolan"
$d.Paragraphs(89).Range.Text = "Denis Villeneuve"

# Block 16: Quel film de science-fiction réalisé par Stanley Kubrick est célèbre pour son intelligence artificielle HAL 9000 ?
$d.Paragraphs(91).Range.Text = "Quel film de science-fiction réalisé par Stanley Kubrick est célèbre pour son intelligence artificielle HAL 9000 ?"
$d.Paragraphs(92).Range.Text = "Blade Runner"
$d.Paragraphs(93).Range.Text = "Alien"
$d.Paragraphs(94).Range.Text = "2001, l'Odyssée de l'espace"
$d.Paragraphs(95).Range.Text = "Matrix"

# Block 17: Quelle actrice a joué le rôle de Cléopâtre dans le film "Cléopâtre" de 1963 ?
$d.Paragraphs(97).Range.Text = "Quelle actrice a joué le rôle de Cléopâtre dans le film `"Cléopâtre`" de 1963 ?"
$d.Paragraphs(98).Range.Text = "Sophia Loren"
$d.Paragraphs(99).Range.Text = "Audrey Hepburn"
$d.Paragraphs(100).Range.Text = "Elizabeth Taylor"
$d.Paragraphs(101).Range.Text = "Ava Gardner"

# Block 18: Quel acteur a joué le rôle de Jack Dawson dans "Titanic" ?
$d.Paragraphs(103).Range.Text = "Quel acteur a joué le rôle de Jack Dawson dans `"Titanic`" ?"
$d.Paragraphs(104).Range.Text = "Leonardo DiCaprio"
$d.Paragraphs(105).Range.Text = "Brad Pitt"
$d.Paragraphs(106).Range.Text = "Tom Hanks"
$d.Paragraphs(107).Range.Text = "Johnny Depp"

# Block 19: Quel film a remporté l'Oscar du meilleur film en 2020 ?
$d.Paragraphs(109).Range.Text = "Quel film a remporté l'Oscar du meilleur film en 2020 ?"
$d.Paragraphs(110).Range.Text = "`"1917`""
$d.Paragraphs(111).Range.Text = "`"The Irishman`""
$d.Paragraphs(112).Range.Text = "`"Parasite`""
$d.Paragraphs(113).Range.Text = "`"Joker`""

# Block 20: Quel est le nom du vaisseau spatial de Han Solo dans "Star Wars" ?
$d.Paragraphs(115).Range.Text = "Quel est le nom du vaisseau spatial de Han Solo dans `"Star Wars`" ?"
$d.Paragraphs(116).Range.Text = "Le Faucon Millenium"
$d.Paragraphs(117).Range.Text = "L'Étoile Noire"
$d.Paragraphs(118).Range.Text = "Le X-Wing"
$d.Paragraphs(119).Range.Text = "Le TIE Fighter"
